$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value before shifting columns
$ws.Range("B2").Value = 16.46212132459583

# Delete entire column C; this shifts D->C and E->D
$ws.Range("C1:C2").EntireColumn.Delete()
